$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Course Materials" sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCM = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCM.Name = "Course Materials"

$wsCode = $wb.Worksheets.Item("Code")

# ---------------------------------------------------------------------------
# 2. Populate "Course Materials" + update "Code" header in the exact order
#    the strings were authored, so the shared-string table comes out in the
#    same order as the target workbook.
# ---------------------------------------------------------------------------
$wsCM.Range("B4").Value = "L001 Course introduction Video "
$wsCM.Range("C1").Value = "Ideation"
$wsCM.Range("E1").Value = "Closing date"
$wsCode.Range("E1").Value = "Closing Date"
$wsCM.Range("D1").Value = "Category"
$wsCM.Range("B2").Value = "Bundle startup (slide templete, tools, folder structure)"
$wsCM.Range("B3").Value = "L001 Course introduction Slides "

$wsCM.Range("F1").Value = "Status"

$wsCM.Range("A2").Value = 1
$wsCM.Range("C2").Value = 43639
$wsCM.Range("E2").Value = 43639
$wsCM.Range("F2").Value = "Done"

$wsCM.Range("A3").Value = 2
$wsCM.Range("C3").Value = 43639
$wsCM.Range("E3").Value = 43639
$wsCM.Range("F3").Value = "Done"

$wsCM.Range("A4").Value = 3
$wsCM.Range("C4").Value = 43639
$wsCM.Range("F4").Value = "ToDo"

# date formatting for the date columns
$wsCM.Range("C2:C3").NumberFormat = "m/d/yyyy"
$wsCM.Range("E2:E3").NumberFormat = "m/d/yyyy"
$wsCM.Range("C4").NumberFormat = "m/d/yyyy"

# header row bold style (matches the other sheets' header style)
$wsCM.Range("A1:F1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Column widths.
# ---------------------------------------------------------------------------
$wsCode.Columns.Item(3).ColumnWidth = 12.5

$wsCM.Columns.Item(1).ColumnWidth = 3.6666666666666665
$wsCM.Columns.Item(2).ColumnWidth = 43.166666666666664
$wsCM.Columns.Item(3).ColumnWidth = 16.5
$wsCM.Columns.Item(4).ColumnWidth = 13.166666666666666
$wsCM.Columns.Item(5).ColumnWidth = 12.666666666666666

# ---------------------------------------------------------------------------
# 4. Selections.
# ---------------------------------------------------------------------------
$wsCM.Activate()
$wsCM.Range("B8").Select()

$wsCode.Activate()
$wsCode.Range("E3").Select()
